$wb = $excel.ActiveWorkbook

# Update "Extracted Date" column (P2:P6) on the Opportunities sheet
$wsOpp = $wb.Worksheets.Item("Opportunities")
$wsOpp.Range("P2").Value = "2025-08-14 03:48"
$wsOpp.Range("P3").Value = "2025-08-14 03:48"
$wsOpp.Range("P4").Value = "2025-08-14 03:48"
$wsOpp.Range("P5").Value = "2025-08-14 03:48"
$wsOpp.Range("P6").Value = "2025-08-14 03:48"

# Update "Generated Date" on the Summary sheet
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B4").Value = "2025-08-14 03:48:11"
